# Update the symbol list (price / volume columns) on the active worksheet.
# Values are stored as literal text in the source data (inline strings),
# so each numeric-looking entry is entered with a leading apostrophe to
# keep Excel from re-typing it as a Number, and the cell style is reset
# back to "Normal" afterwards so no incidental formatting (e.g. a
# quote-prefix / text number format) sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

# --- Column D (Price) updates ---
Set-TextValue "D2"  "245.35"
Set-TextValue "D4"  "5.278"
Set-TextValue "D5"  "0.05773"
Set-TextValue "D6"  "6.458"
Set-TextValue "D7"  "3.143"
Set-TextValue "D8"  "0.8168"
Set-TextValue "D9"  "0.8522"
Set-TextValue "D11" "0.06940"
Set-TextValue "D12" "0.03126"
Set-TextValue "D13" "0.02897"
Set-TextValue "D14" "0.09380"
Set-TextValue "D15" "3.745"
Set-TextValue "D17" "0.04690"
Set-TextValue "D18" "0.0005989"
Set-TextValue "D19" "0.006217"
Set-TextValue "D21" "0.004616"
Set-TextValue "D22" "0.00006899"
Set-TextValue "D23" "3.500"
Set-TextValue "D26" "0.1302"
Set-TextValue "D40" "0.03644"
Set-TextValue "D41" "0.006249"
Set-TextValue "D42" "0.1053"
Set-TextValue "D44" "0.008422"
Set-TextValue "D45" "0.00005270"
Set-TextValue "D47" "0.3699"
Set-TextValue "D48" "0.002283"

# --- Column E (Volume(1h)) updates ---
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
